$wb = $excel.ActiveWorkbook

# --- Fix scenario-name capitalization / formatting on the "schedule" sheet ---
# The header row (row 1) relabels the scenario columns: "low" -> "Low",
# "mid-SC" -> "Mid (SC)", "mid-CC" -> "Mid (CC)", "high" -> "High", and also
# straightens out a swapped pair of columns (B/C) so the "-annual" column
# comes first, matching the rest of the table's column order.
$schedule = $wb.Worksheets.Item("schedule")

$schedule.Range("B1").Value = "Baseline-Low-annual"
$schedule.Range("C1").Value = "Baseline-Low"
$schedule.Range("D1").Value = "Baseline-Mid (SC)-annual"
$schedule.Range("E1").Value = "Baseline-Mid (SC)"
$schedule.Range("F1").Value = "Baseline-Mid (CC)-annual"
$schedule.Range("G1").Value = "Baseline-Mid (CC)"
$schedule.Range("H1").Value = "Moderate-Low-annual"
$schedule.Range("I1").Value = "Moderate-Low"
$schedule.Range("J1").Value = "Moderate-Mid (SC)-annual"
$schedule.Range("K1").Value = "Moderate-Mid (SC)"
$schedule.Range("L1").Value = "Expanded-High-annual"
$schedule.Range("M1").Value = "Expanded-High"

# --- Update the active selection on the "schedule" sheet ---
$schedule.Activate()
$schedule.Range("N1").Select()
